$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter the new donor's figures (replacing the previous organisation's data)
$ws.Range("C9").Value = 139
$ws.Range("C10").Value = 0
$ws.Range("C11").Value = 0
$ws.Range("C12").Value = 101.14

# C9 is a round-euro figure, so it is formatted without decimals
$ws.Range("C9").NumberFormat = "#,##0 ""€"";[Red]-#,##0 ""€"""

# Restore the cursor to the cell the author left selected
[void]$ws.Range("D22").Select()
